# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) for a handful of leve rows across the per-job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 7481.909
$ws.Range("J100").Value = 10186.462
$ws.Range("L100").Value = 10186.462
$ws.Range("N100").Value = -11268.462
$ws.Range("H132").Value = 5310
$ws.Range("I132").Value = 3146.2
$ws.Range("K132").Value = 9438.599999999999
$ws.Range("M132").Value = -6908.599999999999
$ws.Range("H135").Value = 2522.818
$ws.Range("I135").Value = 851.875
$ws.Range("K135").Value = 7666.875
$ws.Range("M135").Value = -5131.875
$ws.Range("H138").Value = 2505.2083
$ws.Range("I138").Value = 2552.15
$ws.Range("J138").Value = 2270.5
$ws.Range("K138").Value = 7656.450000000001
$ws.Range("L138").Value = 6811.5
$ws.Range("M138").Value = -2516.450000000001
$ws.Range("N138").Value = -17091.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4272.1206
$ws.Range("I32").Value = 3760.0566
$ws.Range("K32").Value = 3760.0566
$ws.Range("M32").Value = -3473.0566
$ws.Range("H61").Value = 27503448
$ws.Range("I61").Value = 35004264
$ws.Range("J61").Value = 5000998.5
$ws.Range("K61").Value = 35004264
$ws.Range("L61").Value = 5000998.5
$ws.Range("M61").Value = -35004052
$ws.Range("N61").Value = -5001422.5
$ws.Range("H74").Value = 2595.8965
$ws.Range("J74").Value = 4998.5
$ws.Range("L74").Value = 4998.5
$ws.Range("N74").Value = -6746.5
$ws.Range("H77").Value = 2595.8965
$ws.Range("J77").Value = 4998.5
$ws.Range("L77").Value = 24992.5
$ws.Range("N77").Value = -33728.5
$ws.Range("H132").Value = 3230472.2
$ws.Range("I132").Value = 4989.16
$ws.Range("J132").Value = 16669985
$ws.Range("K132").Value = 14967.48
$ws.Range("L132").Value = 50009955
$ws.Range("M132").Value = -12437.48
$ws.Range("N132").Value = -50015015
$ws.Range("H136").Value = 27503448
$ws.Range("I136").Value = 35004264
$ws.Range("J136").Value = 5000998.5
$ws.Range("K136").Value = 105012792
$ws.Range("L136").Value = 15002995.5
$ws.Range("M136").Value = -105010242
$ws.Range("N136").Value = -15008095.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2252.2307
$ws.Range("I20").Value = 2377.1
$ws.Range("J20").Value = 1836
$ws.Range("K20").Value = 2377.1
$ws.Range("L20").Value = 1836
$ws.Range("M20").Value = -2130.1
$ws.Range("N20").Value = -2330
$ws.Range("H94").Value = 1735.9
$ws.Range("I94").Value = 1984.9
$ws.Range("K94").Value = 1984.9
$ws.Range("M94").Value = -1533.9
$ws.Range("H134").Value = 16670256
$ws.Range("I134").Value = 2880.5
$ws.Range("K134").Value = 8641.5
$ws.Range("M134").Value = -6106.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1969.2667
$ws.Range("I22").Value = 1719
$ws.Range("K22").Value = 1719
$ws.Range("M22").Value = -1369
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H107").Value = 2893.2856
$ws.Range("I107").Value = 1708.8334
$ws.Range("K107").Value = 1708.8334
$ws.Range("M107").Value = 211.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 15963.286
$ws.Range("I14").Value = 15963.286
$ws.Range("K14").Value = 47889.858
$ws.Range("M14").Value = -47716.858
$ws.Range("H51").Value = 11310.667
$ws.Range("I51").Value = 299.5
$ws.Range("K51").Value = 898.5
$ws.Range("M51").Value = -438.5
$ws.Range("H68").Value = 4124.75
$ws.Range("J68").Value = 4666.3335
$ws.Range("L68").Value = 13999.0005
$ws.Range("N68").Value = -15621.0005
$ws.Range("H71").Value = 4124.75
$ws.Range("J71").Value = 4666.3335
$ws.Range("L71").Value = 41997.0015
$ws.Range("N71").Value = -50109.0015
$ws.Range("H86").Value = 2607.6667
$ws.Range("J86").Value = 3879.5
$ws.Range("L86").Value = 11638.5
$ws.Range("N86").Value = -14010.5
$ws.Range("H89").Value = 2607.6667
$ws.Range("J89").Value = 3879.5
$ws.Range("L89").Value = 34915.5
$ws.Range("N89").Value = -46771.5
$ws.Range("H112").Value = 19072
$ws.Range("I112").Value = 12027
$ws.Range("J112").Value = 20833.25
$ws.Range("K112").Value = 36081
$ws.Range("L112").Value = 62499.75
$ws.Range("M112").Value = -34973
$ws.Range("N112").Value = -64715.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7939.3823
$ws.Range("I70").Value = 7452.1177
$ws.Range("J70").Value = 8426.647000000001
$ws.Range("K70").Value = 7452.1177
$ws.Range("L70").Value = 8426.647000000001
$ws.Range("M70").Value = -7182.1177
$ws.Range("N70").Value = -8966.647000000001
$ws.Range("H73").Value = 7939.3823
$ws.Range("I73").Value = 7452.1177
$ws.Range("J73").Value = 8426.647000000001
$ws.Range("K73").Value = 7452.1177
$ws.Range("L73").Value = 8426.647000000001
$ws.Range("M73").Value = -6516.1177
$ws.Range("N73").Value = -10298.647
$ws.Range("H102").Value = 2827.3572
$ws.Range("I102").Value = 2660.3076
$ws.Range("K102").Value = 2660.3076
$ws.Range("M102").Value = -1038.3076
$ws.Range("H113").Value = 1545330.4
$ws.Range("I113").Value = 2161.5557
$ws.Range("J113").Value = 6174836.5
$ws.Range("K113").Value = 2161.5557
$ws.Range("L113").Value = 6174836.5
$ws.Range("M113").Value = 8.444300000000112
$ws.Range("N113").Value = -6179176.5
$ws.Range("H132").Value = 12503911
$ws.Range("I132").Value = 2822.25
$ws.Range("K132").Value = 8466.75
$ws.Range("M132").Value = -5936.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 99999
$ws.Range("J110").Value = 99999
$ws.Range("L110").Value = 99999
$ws.Range("N110").Value = -108179
$ws.Range("H122").Value = 3107.8914
$ws.Range("I122").Value = 2961.0715
$ws.Range("J122").Value = 4649.5
$ws.Range("K122").Value = 8883.2145
$ws.Range("L122").Value = 13948.5
$ws.Range("M122").Value = -6433.2145
$ws.Range("N122").Value = -18848.5
$ws.Range("H132").Value = 4217.615
$ws.Range("I132").Value = 2504.125
$ws.Range("K132").Value = 7512.375
$ws.Range("M132").Value = -4982.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14654.5
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 14654.5
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 14654.5
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -15434.5
$ws.Range("H132").Value = 719760.6
$ws.Range("I132").Value = 4765.4
$ws.Range("K132").Value = 14296.2
$ws.Range("M132").Value = -11766.2
